$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 15
$ws.Range("H15").Value = 1132.7567
$ws.Range("I15").Value = 1132.7567
$ws.Range("K15").Value = 3398.2701
$ws.Range("M15").Value = -3229.2701

# row 17
$ws.Range("H17").Value = 1519.3549
$ws.Range("J17").Value = 1519.3549
$ws.Range("L17").Value = 4558.0647
$ws.Range("N17").Value = -4894.0647

# row 113
$ws.Range("H113").Value = 5200
$ws.Range("I113").Value = 6280.4
$ws.Range("K113").Value = 6280.4
$ws.Range("M113").Value = -3026.4

# row 129
$ws.Range("H129").Value = 1724.5
$ws.Range("I129").Value = 1724.5
$ws.Range("K129").Value = 5173.5
$ws.Range("M129").Value = -173.5

# row 132
$ws.Range("H132").Value = 3139.7827
$ws.Range("I132").Value = 2831.75
$ws.Range("K132").Value = 8495.25
$ws.Range("M132").Value = -5965.25

# row 137
$ws.Range("H137").Value = 1781.1482
$ws.Range("I137").Value = 1601.95
$ws.Range("J137").Value = 2293.1428
$ws.Range("K137").Value = 4805.85
$ws.Range("L137").Value = 6879.428400000001
$ws.Range("M137").Value = -2255.85
$ws.Range("N137").Value = -11979.4284

# row 138
$ws.Range("H138").Value = 3148.8914
$ws.Range("J138").Value = 3566.121
$ws.Range("L138").Value = 10698.363
$ws.Range("N138").Value = -20978.363


$ws = $wb.Worksheets.Item("ARM")
# row 2
$ws.Range("H2").Value = 1596.4
$ws.Range("I2").Value = 1596.4
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1596.4
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -1483.4
$ws.Range("N2").ClearContents()

# row 32
$ws.Range("H32").Value = 5555.446
$ws.Range("I32").Value = 552.0678
$ws.Range("K32").Value = 552.0678
$ws.Range("M32").Value = -265.0678

# row 45
$ws.Range("H45").Value = 15999
$ws.Range("I45").Value = 17892.54
$ws.Range("K45").Value = 17892.54
$ws.Range("M45").Value = -17515.54

# row 116
$ws.Range("H116").Value = 1596.4
$ws.Range("I116").Value = 1596.4
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1596.4
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 697.5999999999999
$ws.Range("N116").ClearContents()

# row 122
$ws.Range("H122").Value = 3899.6667
$ws.Range("I122").Value = 3679.6
$ws.Range("K122").Value = 11038.8
$ws.Range("M122").Value = -8588.799999999999

# row 132
$ws.Range("H132").Value = 2378.68
$ws.Range("I132").Value = 2248.6365
$ws.Range("K132").Value = 6745.9095
$ws.Range("M132").Value = -4215.9095


$ws = $wb.Worksheets.Item("BSM")
# row 3
$ws.Range("H3").Value = 1596.4
$ws.Range("I3").Value = 1596.4
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1596.4
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -1482.4
$ws.Range("N3").ClearContents()

# row 22
$ws.Range("H22").Value = 585.75
$ws.Range("I22").Value = 747
$ws.Range("J22").Value = 102
$ws.Range("K22").Value = 747
$ws.Range("L22").Value = 102
$ws.Range("M22").Value = -574
$ws.Range("N22").Value = -448


$ws = $wb.Worksheets.Item("CRP")
# row 16
$ws.Range("H16").Value = 21934
$ws.Range("I16").Value = 13673.083
$ws.Range("J16").Value = 71499.5
$ws.Range("K16").Value = 13673.083
$ws.Range("L16").Value = 71499.5
$ws.Range("M16").Value = -13386.083
$ws.Range("N16").Value = -72073.5

# row 31
$ws.Range("H31").Value = 4453.0967
$ws.Range("I31").Value = 5338.8125
$ws.Range("J31").Value = 3508.3333
$ws.Range("K31").Value = 5338.8125
$ws.Range("L31").Value = 3508.3333
$ws.Range("M31").Value = -5043.8125
$ws.Range("N31").Value = -4098.3333

# row 34
$ws.Range("H34").Value = 4453.0967
$ws.Range("I34").Value = 5338.8125
$ws.Range("J34").Value = 3508.3333
$ws.Range("K34").Value = 5338.8125
$ws.Range("L34").Value = 3508.3333
$ws.Range("M34").Value = -5136.8125
$ws.Range("N34").Value = -3912.3333

# row 62
$ws.Range("H62").Value = 3653.5557
$ws.Range("I62").Value = 3698
$ws.Range("J62").Value = 3618
$ws.Range("K62").Value = 3698
$ws.Range("L62").Value = 3618
$ws.Range("M62").Value = -3074
$ws.Range("N62").Value = -4866

# row 65
$ws.Range("H65").Value = 3653.5557
$ws.Range("I65").Value = 3698
$ws.Range("J65").Value = 3618
$ws.Range("K65").Value = 18490
$ws.Range("L65").Value = 18090
$ws.Range("M65").Value = -15370
$ws.Range("N65").Value = -24330

# row 99
$ws.Range("H99").Value = 6482.25
$ws.Range("I99").Value = 6898
$ws.Range("J99").Value = 6185.2856
$ws.Range("K99").Value = 6898
$ws.Range("L99").Value = 6185.2856
$ws.Range("M99").Value = -5400
$ws.Range("N99").Value = -9181.285599999999

# row 113
$ws.Range("H113").Value = 21934
$ws.Range("I113").Value = 13673.083
$ws.Range("J113").Value = 71499.5
$ws.Range("K113").Value = 13673.083
$ws.Range("L113").Value = 71499.5
$ws.Range("M113").Value = -11503.083
$ws.Range("N113").Value = -75839.5

# row 126
$ws.Range("H126").Value = 6482.25
$ws.Range("I126").Value = 6898
$ws.Range("J126").Value = 6185.2856
$ws.Range("K126").Value = 20694
$ws.Range("L126").Value = 18555.8568
$ws.Range("M126").Value = -18224
$ws.Range("N126").Value = -23495.8568


$ws = $wb.Worksheets.Item("CUL")
# row 2
$ws.Range("H2").Value = 2038.75
$ws.Range("I2").Value = 2171.2856
$ws.Range("J2").Value = 1111
$ws.Range("K2").Value = 13027.7136
$ws.Range("L2").Value = 6666
$ws.Range("M2").Value = -12914.7136
$ws.Range("N2").Value = -6892

# row 5
$ws.Range("H5").Value = 789.8
$ws.Range("I5").Value = 670.4286
$ws.Range("J5").Value = 1068.3334
$ws.Range("K5").Value = 2011.2858
$ws.Range("L5").Value = 3205.0002
$ws.Range("M5").Value = -1899.2858
$ws.Range("N5").Value = -3429.0002

# row 11
$ws.Range("H11").Value = 2083.1667
$ws.Range("I11").Value = 480
$ws.Range("J11").Value = 3686.3333
$ws.Range("K11").Value = 1440
$ws.Range("L11").Value = 11058.9999
$ws.Range("M11").Value = -1300
$ws.Range("N11").Value = -11338.9999

# row 22
$ws.Range("H22").Value = 1499
$ws.Range("J22").Value = 1499
$ws.Range("L22").Value = 4497
$ws.Range("N22").Value = -4835

# row 26
$ws.Range("H26").Value = 1450
$ws.Range("I26").Value = 1400
$ws.Range("J26").Value = 1500
$ws.Range("K26").Value = 4200
$ws.Range("L26").Value = 4500
$ws.Range("M26").Value = -3912
$ws.Range("N26").Value = -5076

# row 27
$ws.Range("H27").Value = 1499
$ws.Range("J27").Value = 1499
$ws.Range("L27").Value = 4497
$ws.Range("N27").Value = -4701

# row 92
$ws.Range("H92").Value = 455.83334
$ws.Range("I92").Value = 400
$ws.Range("J92").Value = 483.75
$ws.Range("K92").Value = 1200
$ws.Range("L92").Value = 1451.25
$ws.Range("M92").Value = 48
$ws.Range("N92").Value = -3947.25

# row 134
$ws.Range("H134").Value = 630.1429000000001
$ws.Range("I134").Value = 630.1429000000001
$ws.Range("K134").Value = 1890.4287
$ws.Range("M134").Value = 3179.5713

# row 135
$ws.Range("H135").Value = 789.8
$ws.Range("I135").Value = 670.4286
$ws.Range("J135").Value = 1068.3334
$ws.Range("K135").Value = 6033.8574
$ws.Range("L135").Value = 9615.000599999999
$ws.Range("M135").Value = -3498.8574
$ws.Range("N135").Value = -14685.0006


$ws = $wb.Worksheets.Item("GSM")
# row 122
$ws.Range("H122").Value = 1999
$ws.Range("I122").Value = 1998.75
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 5996.25
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -3546.25
$ws.Range("N122").Value = -10900

# row 123
$ws.Range("H123").Value = 34449.832
$ws.Range("J123").Value = 34449.832
$ws.Range("L123").Value = 34449.832
$ws.Range("N123").Value = -39349.832

# row 132
$ws.Range("H132").Value = 3590.0962
$ws.Range("I132").Value = 3355.739
$ws.Range("K132").Value = 10067.217
$ws.Range("M132").Value = -7537.217000000001


$ws = $wb.Worksheets.Item("LTW")
# row 132
$ws.Range("H132").Value = 2300.4
$ws.Range("I132").Value = 2167.6667
$ws.Range("K132").Value = 6503.000100000001
$ws.Range("M132").Value = -3973.000100000001


$ws = $wb.Worksheets.Item("WVR")
# row 122
$ws.Range("H122").Value = 2681
$ws.Range("I122").Value = 2559.7334
$ws.Range("K122").Value = 7679.2002
$ws.Range("M122").Value = -5229.2002

